# Project "Sample Project" save: update rule row R40's Rule label (B11)
# from "R40" to "1" on the Rules sheet. Entering a leading apostrophe
# forces the numeric-looking text to be stored as a text value (shared
# string), matching how Excel stores literal text that looks like a number.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Rules")
if (-not $ws) {
    $ws = $wb.ActiveSheet
}

$ws.Range("B11").Value = "'1"
